$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
  2  = @{ E=3; G=40.42873533333334; H=121.286206;        I=0.08313576592793961; J=0.08313576592793961; K=3; M=16.535604;        N=49.606812;        O=0.2120453146491552; P=0.2120453146491552; Q=668.513557692808;  R=6016.622019235272;  S=0.01762854964478847; T=0.01762854964478848 }
  3  = @{ E=3; G=40.42873533333334; H=121.286206;        I=0.08313576592793961; J=0.08313576592793961; K=3; M=40.62063066666667; N=121.861892;       O=0.5209011059384622; P=0.5209011059384622; Q=1642.24072629575;  R=14780.16653666175;  S=0.04330551241490486; T=0.04330551241490486 }
  4  = @{ E=3; G=40.42873533333334; H=121.286206;        I=0.08313576592793961; J=0.08313576592793961; K=3; M=20.825229;        N=62.475687;        O=0.2670535794123827; P=0.2670535794123827; Q=841.9376714970581; R=7577.439043473522;  S=0.02220170386824628; T=0.02220170386824628 }
  5  = @{ E=3; G=412.4720866666667; H=1237.41626;         I=0.84818836320749;    J=0.84818836320749;    K=3; M=16.535604;        N=49.606812;        O=0.2120453146491552; P=0.2120453146491552; Q=6820.475086173679; R=61384.27577556312;  S=0.1798543683580842;  T=0.1798543683580842  }
  6  = @{ E=3; G=412.4720866666667; H=1237.41626;         I=0.84818836320749;    J=0.84818836320749;    K=3; M=40.62063066666667; N=121.861892;       O=0.5209011059384622; P=0.5209011059384622; Q=16754.87629279599; R=150793.8866351639;  S=0.4418222564389156;  T=0.4418222564389156  }
  7  = @{ E=3; G=412.4720866666667; H=1237.41626;         I=0.84818836320749;    J=0.84818836320749;    K=3; M=20.825229;        N=62.475687;        O=0.2670535794123827; P=0.2670535794123827; Q=8589.82566094118;  R=77308.43094847062;  S=0.2265117384104903;  T=0.2265117384104903  }
  8  = @{ E=3; G=33.396921;         H=100.190763;         I=0.0686758708645703;  J=0.0686758708645703;  K=3; M=16.535604;        N=49.606812;        O=0.2120453146491552; P=0.2120453146491552; Q=552.2382604752839; R=4970.144344277555;  S=0.01456239664628256; T=0.01456239664628256 }
  9  = @{ E=3; G=33.396921;         H=100.190763;         I=0.0686758708645703;  J=0.0686758708645703;  K=3; M=40.62063066666667; N=121.861892;       O=0.5209011059384622; P=0.5209011059384622; Q=1356.603993344844; R=12209.4359401036;   S=0.03577333708464168; T=0.03577333708464168 }
  10 = @{ E=3; G=33.396921;         H=100.190763;         I=0.0686758708645703;  J=0.0686758708645703;  K=3; M=20.825229;        N=62.475687;        O=0.2670535794123827; P=0.2670535794123827; Q=695.498527719909;  R=6259.486749479181;  S=0.01834013713364606; T=0.01834013713364606 }
}

foreach ($row in $data.Keys) {
  $cols = $data[$row]
  foreach ($col in $cols.Keys) {
    $ws.Range("$col$row").Value = $cols[$col]
  }
}
